$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.306.20"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.869.11"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.12"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06575"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.59"
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08017"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.86"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").Value = "1.871.07"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6845"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.68"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("D17").Value = "30.320.43"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.03"
$ws.Range("E18").Value = "  +3.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007609"
$ws.Range("E19").Value = "  +3.57%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "2.116.45"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.263"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.212"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.390"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.63"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.950"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.368"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09874"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.359"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.064"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04711"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.136"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6998"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01880"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.696"
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.260"
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.90"
$ws.Range("E41").Value = "  -6.40%  "
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8418"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4163"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.79"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.053"
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.093"
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "911.93"
$ws.Range("E49").Value = "  -4.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.45"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05703"
$ws.Range("E51").Value = "  +1.15%  "
